$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 33
$ws.Range("H33").Value = 494.3846
$ws.Range("I33").Value = 347.18182
$ws.Range("K33").Value = 347.18182
$ws.Range("M33").Value = -118.18182
# Row 53
$ws.Range("H53").Value = 1205.4706
$ws.Range("I53").Value = 1143.2858
$ws.Range("J53").Value = 1249
$ws.Range("K53").Value = 1143.2858
$ws.Range("L53").Value = 1249
$ws.Range("M53").Value = -506.2858000000001
$ws.Range("N53").Value = -2523
# Row 64
$ws.Range("H64").Value = 6896.385
$ws.Range("I64").Value = 4664.5713
$ws.Range("K64").Value = 4664.5713
$ws.Range("M64").Value = -4416.5713
# Row 67
$ws.Range("H67").Value = 6896.385
$ws.Range("I67").Value = 4664.5713
$ws.Range("K67").Value = 4664.5713
$ws.Range("M67").Value = -3806.5713
# Row 74
$ws.Range("I74").Value = 9162.625
$ws.Range("J74").Value = 7166.6665
$ws.Range("K74").Value = 9162.625
$ws.Range("L74").Value = 7166.6665
$ws.Range("M74").Value = -8226.625
$ws.Range("N74").Value = -9038.666499999999
# Row 77
$ws.Range("I77").Value = 9162.625
$ws.Range("J77").Value = 7166.6665
$ws.Range("K77").Value = 45813.125
$ws.Range("L77").Value = 35833.3325
$ws.Range("M77").Value = -41133.125
$ws.Range("N77").Value = -45193.3325
# Row 115
$ws.Range("H115").Value = 909.7143
$ws.Range("I115").Value = 909.7143
$ws.Range("J115").Value = 0
$ws.Range("K115").Value = 2729.1429
$ws.Range("L115").Value = 0
$ws.Range("M115").Value = -1162.1429
$ws.Range("N115").ClearContents()
# Row 138
$ws.Range("H138").Value = 1737.46
$ws.Range("I138").Value = 796.6047
$ws.Range("J138").Value = 2447.228
$ws.Range("K138").Value = 2389.8141
$ws.Range("L138").Value = 7341.684
$ws.Range("M138").Value = 2750.1859
$ws.Range("N138").Value = -17621.684

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 5
$ws.Range("H5").Value = 676.4167
$ws.Range("I5").Value = 819
$ws.Range("J5").Value = 248.66667
$ws.Range("K5").Value = 819
$ws.Range("L5").Value = 248.66667
$ws.Range("M5").Value = -707
$ws.Range("N5").Value = -472.66667
# Row 32
$ws.Range("H32").Value = 21075498
$ws.Range("I32").Value = 24986588
$ws.Range("K32").Value = 24986588
$ws.Range("M32").Value = -24986301
# Row 74
$ws.Range("H74").Value = 2250.2954
$ws.Range("I74").Value = 2070.0698
$ws.Range("K74").Value = 2070.0698
$ws.Range("M74").Value = -1196.0698
# Row 77
$ws.Range("H77").Value = 2250.2954
$ws.Range("I77").Value = 2070.0698
$ws.Range("K77").Value = 10350.349
$ws.Range("M77").Value = -5982.349000000002

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 4
$ws.Range("H4").Value = 635.8461
$ws.Range("I4").Value = 752
$ws.Range("J4").Value = 248.66667
$ws.Range("K4").Value = 752
$ws.Range("L4").Value = 248.66667
$ws.Range("M4").Value = -637
$ws.Range("N4").Value = -478.66667
# Row 107
$ws.Range("H107").Value = 2420.7273
$ws.Range("I107").Value = 2028.625
$ws.Range("J107").Value = 3466.3333
$ws.Range("K107").Value = 2028.625
$ws.Range("L107").Value = 3466.3333
$ws.Range("M107").Value = -108.625
$ws.Range("N107").Value = -7306.3333
# Row 131
$ws.Range("H131").Value = 112494
$ws.Range("J131").Value = 112494
$ws.Range("L131").Value = 112494
$ws.Range("N131").Value = -122574
# Row 134
$ws.Range("H134").Value = 1787973.9
$ws.Range("I134").Value = 2102060.5
$ws.Range("J134").Value = 8150
$ws.Range("K134").Value = 6306181.5
$ws.Range("L134").Value = 24450
$ws.Range("M134").Value = -6303646.5
$ws.Range("N134").Value = -29520

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 58
$ws.Range("H58").Value = 2585.6
$ws.Range("I58").Value = 2508.7568
$ws.Range("J58").Value = 3533.3333
$ws.Range("K58").Value = 2508.7568
$ws.Range("L58").Value = 3533.3333
$ws.Range("M58").Value = -2305.7568
$ws.Range("N58").Value = -3939.3333
# Row 122
$ws.Range("H122").Value = 2945527.5
$ws.Range("I122").Value = 4351683
$ws.Range("K122").Value = 13055049
$ws.Range("M122").Value = -13052599
# Row 136
$ws.Range("H136").Value = 2585.6
$ws.Range("I136").Value = 2508.7568
$ws.Range("J136").Value = 3533.3333
$ws.Range("K136").Value = 7526.2704
$ws.Range("L136").Value = 10599.9999
$ws.Range("M136").Value = -4976.2704
$ws.Range("N136").Value = -15699.9999

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 122
$ws.Range("H122").Value = 637499.5600000001
$ws.Range("J122").Value = 871701.9
$ws.Range("L122").Value = 7845317.100000001
$ws.Range("N122").Value = -7850217.100000001
# Row 140
$ws.Range("H140").Value = 2093.9167
$ws.Range("I140").Value = 1291.8889
$ws.Range("K140").Value = 3875.6667
$ws.Range("M140").Value = 1304.3333

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 22151.895
$ws.Range("I70").Value = 67077.39999999999
$ws.Range("J70").Value = 6107.0713
$ws.Range("K70").Value = 67077.39999999999
$ws.Range("L70").Value = 6107.0713
$ws.Range("M70").Value = -66807.39999999999
$ws.Range("N70").Value = -6647.0713
# Row 73
$ws.Range("H73").Value = 22151.895
$ws.Range("I73").Value = 67077.39999999999
$ws.Range("J73").Value = 6107.0713
$ws.Range("K73").Value = 67077.39999999999
$ws.Range("L73").Value = 6107.0713
$ws.Range("M73").Value = -66141.39999999999
$ws.Range("N73").Value = -7979.0713
# Row 126
$ws.Range("H126").Value = 2912.7778
$ws.Range("I126").Value = 2798.75
$ws.Range("J126").Value = 3004
$ws.Range("K126").Value = 8396.25
$ws.Range("L126").Value = 9012
$ws.Range("M126").Value = -5926.25
$ws.Range("N126").Value = -13952

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 5306.6665
$ws.Range("I7").Value = 5742.222
$ws.Range("K7").Value = 5742.222
$ws.Range("M7").Value = -5630.222
# Row 22
$ws.Range("H22").Value = 2580.6667
$ws.Range("I22").Value = 1104.3334
$ws.Range("J22").Value = 5533.3335
$ws.Range("K22").Value = 1104.3334
$ws.Range("L22").Value = 5533.3335
$ws.Range("M22").Value = -809.3334
$ws.Range("N22").Value = -6123.3335
# Row 27
$ws.Range("H27").Value = 2580.6667
$ws.Range("I27").Value = 1104.3334
$ws.Range("J27").Value = 5533.3335
$ws.Range("K27").Value = 1104.3334
$ws.Range("L27").Value = 5533.3335
$ws.Range("M27").Value = -997.3334
$ws.Range("N27").Value = -5747.3335
# Row 40
$ws.Range("H40").Value = 37040144
$ws.Range("I40").Value = 41668290
$ws.Range("K40").Value = 41668290
$ws.Range("M40").Value = -41668154
# Row 46
$ws.Range("H46").Value = 3693.4
$ws.Range("J46").Value = 4474.7334
$ws.Range("L46").Value = 4474.7334
$ws.Range("N46").Value = -4850.7334
# Row 55
$ws.Range("H55").Value = 215.55556
$ws.Range("I55").Value = 162
$ws.Range("K55").Value = 162
$ws.Range("M55").Value = 11
# Row 122
$ws.Range("H122").Value = 8568.272000000001
$ws.Range("I122").Value = 8211.277
$ws.Range("J122").Value = 10174.75
$ws.Range("K122").Value = 24633.831
$ws.Range("L122").Value = 30524.25
$ws.Range("M122").Value = -22183.831
$ws.Range("N122").Value = -35424.25
# Row 126
$ws.Range("H126").Value = 5306.6665
$ws.Range("I126").Value = 5742.222
$ws.Range("K126").Value = 17226.666
$ws.Range("M126").Value = -14756.666
# Row 136
$ws.Range("H136").Value = 1840.6129
$ws.Range("J136").Value = 1676
$ws.Range("L136").Value = 5028
$ws.Range("N136").Value = -10128

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 81
$ws.Range("H81").Value = 81489.46000000001
$ws.Range("I81").Value = 170244.17
$ws.Range("K81").Value = 340488.34
$ws.Range("M81").Value = -339427.34
# Row 84
$ws.Range("H84").Value = 81489.46000000001
$ws.Range("I84").Value = 170244.17
$ws.Range("K84").Value = 1702441.7
$ws.Range("M84").Value = -1697137.7
# Row 113
$ws.Range("H113").Value = 562.93335
$ws.Range("I113").Value = 516.1667
$ws.Range("J113").Value = 750
$ws.Range("K113").Value = 1548.5001
$ws.Range("L113").Value = 2250
$ws.Range("M113").Value = 621.4999
$ws.Range("N113").Value = -6590
# Row 122
$ws.Range("H122").Value = 35719116
$ws.Range("I122").Value = 43483536
$ws.Range("J122").Value = 2799.8
$ws.Range("K122").Value = 130450608
$ws.Range("L122").Value = 8399.400000000001
$ws.Range("M122").Value = -130448158
$ws.Range("N122").Value = -13299.4
# Row 124
$ws.Range("H124").Value = 65992.336
$ws.Range("J124").Value = 65992.336
$ws.Range("L124").Value = 65992.336
$ws.Range("N124").Value = -75812.336
# Row 125
$ws.Range("H125").Value = 49994.5
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 49994.5
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 49994.5
$ws.Range("N125").Value = -59834.5
$ws.Range("M125").ClearContents()
# Row 132
$ws.Range("H132").Value = 30041.027
$ws.Range("I132").Value = 33155.22
$ws.Range("K132").Value = 99465.66
$ws.Range("M132").Value = -96935.66

